# Apply "Upload new analysis run" changes to the keyword co-occurrence
# table on the active worksheet (rows 2-46, columns B/C/D).
#  - Column B/C keyword labels are updated: some are renamed (e.g.
#    "Energy Transition" -> "green transition", "Biodiversity" ->
#    "loss of biodiversity", "Glaciers" -> "melting glaciers",
#    "Fake News" -> "misinformation") and most others are simply
#    lower-cased to match the new run's formatting.
#  - Column D co-occurrence counts are refreshed with the new run's
#    numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "green transition"
$ws.Range("C2").Value = "greenhouse effect"
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = "green transition"
$ws.Range("C3").Value = "loss of biodiversity"
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = "green transition"
$ws.Range("C4").Value = "extreme weather events"
$ws.Range("D4").Value = 0

$ws.Range("B5").Value = "green transition"
$ws.Range("C5").Value = "CO2"
$ws.Range("D5").Value = 2

$ws.Range("B6").Value = "green transition"
$ws.Range("C6").Value = "emissions"
$ws.Range("D6").Value = 7

$ws.Range("B7").Value = "green transition"
$ws.Range("C7").Value = "global warming"
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = "green transition"
$ws.Range("C8").Value = "melting glaciers"
$ws.Range("D8").Value = 0

$ws.Range("B9").Value = "green transition"
$ws.Range("C9").Value = "renewable energy"
$ws.Range("D9").Value = 6

$ws.Range("B10").Value = "green transition"
$ws.Range("C10").Value = "misinformation"
$ws.Range("D10").Value = 0

$ws.Range("B11").Value = "greenhouse effect"
$ws.Range("C11").Value = "loss of biodiversity"
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = "greenhouse effect"
$ws.Range("C12").Value = "extreme weather events"
$ws.Range("D12").Value = 1

$ws.Range("B13").Value = "greenhouse effect"
$ws.Range("C13").Value = "CO2"
$ws.Range("D13").Value = 80

$ws.Range("B14").Value = "greenhouse effect"
$ws.Range("C14").Value = "emissions"
$ws.Range("D14").Value = 32

$ws.Range("B15").Value = "greenhouse effect"
$ws.Range("C15").Value = "global warming"
$ws.Range("D15").Value = 41

$ws.Range("B16").Value = "greenhouse effect"
$ws.Range("C16").Value = "melting glaciers"
$ws.Range("D16").Value = 0

$ws.Range("B17").Value = "greenhouse effect"
$ws.Range("C17").Value = "renewable energy"
$ws.Range("D17").Value = 1

$ws.Range("B18").Value = "greenhouse effect"
$ws.Range("C18").Value = "misinformation"
$ws.Range("D18").Value = 0

$ws.Range("B19").Value = "loss of biodiversity"
$ws.Range("C19").Value = "extreme weather events"
$ws.Range("D19").Value = 1

$ws.Range("B20").Value = "loss of biodiversity"
$ws.Range("C20").Value = "CO2"
$ws.Range("D20").Value = 1

$ws.Range("B21").Value = "loss of biodiversity"
$ws.Range("C21").Value = "emissions"
$ws.Range("D21").Value = 7

$ws.Range("B22").Value = "loss of biodiversity"
$ws.Range("C22").Value = "global warming"
$ws.Range("D22").Value = 2

$ws.Range("B23").Value = "loss of biodiversity"
$ws.Range("C23").Value = "melting glaciers"
$ws.Range("D23").Value = 1

$ws.Range("B24").Value = "loss of biodiversity"
$ws.Range("C24").Value = "renewable energy"
$ws.Range("D24").Value = 0

$ws.Range("B25").Value = "loss of biodiversity"
$ws.Range("C25").Value = "misinformation"
$ws.Range("D25").Value = 0

$ws.Range("B26").Value = "extreme weather events"
$ws.Range("C26").Value = "CO2"
$ws.Range("D26").Value = 2

$ws.Range("B27").Value = "extreme weather events"
$ws.Range("C27").Value = "emissions"
$ws.Range("D27").Value = 16

$ws.Range("B28").Value = "extreme weather events"
$ws.Range("C28").Value = "global warming"
$ws.Range("D28").Value = 23

$ws.Range("B29").Value = "extreme weather events"
$ws.Range("C29").Value = "melting glaciers"
$ws.Range("D29").Value = 4

$ws.Range("B30").Value = "extreme weather events"
$ws.Range("C30").Value = "renewable energy"
$ws.Range("D30").Value = 3

$ws.Range("B31").Value = "extreme weather events"
$ws.Range("C31").Value = "misinformation"
$ws.Range("D31").Value = 4

$ws.Range("B32").Value = "CO2"
$ws.Range("C32").Value = "emissions"
$ws.Range("D32").Value = 12

$ws.Range("B33").Value = "CO2"
$ws.Range("C33").Value = "global warming"
$ws.Range("D33").Value = 4

$ws.Range("B34").Value = "CO2"
$ws.Range("C34").Value = "melting glaciers"
$ws.Range("D34").Value = 0

$ws.Range("B35").Value = "CO2"
$ws.Range("C35").Value = "renewable energy"
$ws.Range("D35").Value = 0

$ws.Range("B36").Value = "CO2"
$ws.Range("C36").Value = "misinformation"
$ws.Range("D36").Value = 0

$ws.Range("B37").Value = "emissions"
$ws.Range("C37").Value = "global warming"
$ws.Range("D37").Value = 19

$ws.Range("B38").Value = "emissions"
$ws.Range("C38").Value = "melting glaciers"
$ws.Range("D38").Value = 6

$ws.Range("B39").Value = "emissions"
$ws.Range("C39").Value = "renewable energy"
$ws.Range("D39").Value = 1

$ws.Range("B40").Value = "emissions"
$ws.Range("C40").Value = "misinformation"
$ws.Range("D40").Value = 0

$ws.Range("B41").Value = "global warming"
$ws.Range("C41").Value = "melting glaciers"
$ws.Range("D41").Value = 24

$ws.Range("B42").Value = "global warming"
$ws.Range("C42").Value = "renewable energy"
$ws.Range("D42").Value = 1

$ws.Range("B43").Value = "global warming"
$ws.Range("C43").Value = "misinformation"
$ws.Range("D43").Value = 0

$ws.Range("B44").Value = "melting glaciers"
$ws.Range("C44").Value = "renewable energy"
$ws.Range("D44").Value = 2

$ws.Range("B45").Value = "melting glaciers"
$ws.Range("C45").Value = "misinformation"
$ws.Range("D45").Value = 0

$ws.Range("B46").Value = "renewable energy"
$ws.Range("C46").Value = "misinformation"
$ws.Range("D46").Value = 0
